# Commit: "fixed xls files execution"
#
# The sheet name contained a plain space ("Контроль Форма22"), which trips up
# downstream automation that parses/loads the sheet name without quoting it
# (e.g. building A1-style or file-based references). Replace the space with
# an underscore so the name is a single "word" token: "Контроль_Форма22".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Контроль_Форма22"
